$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sponsor activity toggle cell: "???" -> "NULL"
$ws.Range("A2").Value = "NULL"

# Move the active selection to A3 (matches the recorded cursor position)
$ws.Range("A3").Select()
